# Replace the poem texts in A2:A11, resize each row to match the new
# content, and update the active selection (matches the upstream xlsx diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Cemal Süreya - "Desem ki":
Desem ki vakitlerden bir nisan akşamıdır,
Rüzgârların en ferahlatıcısı senden esiyor,
Sende seyrediyorum denizlerin en mavisini,
Ormanların en kuytusunu sende gezmekteyim,
Senden kopardım çiçeklerin en solmazını,
Toprakların en bereketlisini sende sürdüm,
Sende tattım yemişlerin cümlesini.'
$ws.Rows.Item(2).RowHeight = 185.4

$ws.Range("A3").Value = 'Nazım Hikmet - "Seviyorum Seni":
Seviyorum seni ekmeği tuza banıp yer gibi
geceleyin ateşler içinde uyanarak
ağzımı dayayıp musluğa su içer gibi,
ağır posta paketini, neyin nesi belirsiz,
telaşlı, sevinçli, kuşkulu açar gibi,
seviyorum seni denizi uçakla ilk defa geçer gibi. '
$ws.Rows.Item(3).RowHeight = 145.8

$ws.Range("A4").Value = 'Özdemir Asaf - "Lavinia":
Sana gitme demeyeceğim.
Üşüyorsun ceketimi al.
Günün en güzel saatleri bunlar.
Yanımda kal.
Sana gitme demeyeceğim.
Ama gitme, Lavinia. '
$ws.Rows.Item(4).RowHeight = 93

$ws.Range("A5").Value = 'Attila İlhan - "Bir Adın Kalmalı":
Bir Adın Kalmalı
Geriye
Bütün kırılmış şeylerin
Nihayetinde
Aynaların ardında sır
Yalnızlığın peşinde kuvvet
Evet
Nihayet
Bir Adın Kalmalı
Geriye. '
$ws.Rows.Item(5).RowHeight = 145.8

$ws.Range("A6").Value = 'Turgut Uyar - "Hasret Bir Şey Değil":
Hasret bir şey değil,
Elagözlüm
Ömrümüz böyle olmamalıydı
Hep aşkta durmalıydı çağımız.'
$ws.Rows.Item(6).RowHeight = 79.8

$ws.Range("A7").Value = 'Ahmed Arif - "Hasretinden Prangalar Eskittim" 
Seni anlatabilmek seni.
İyi çocuklara, kahramanlara.
Seni anlatabilmek seni,
Namussuz bir dünya sofrasında,
Tırnaklarımda et, dişlerimde kan…'
$ws.Rows.Item(7).RowHeight = 100.8

$ws.Range("A8").Value = 'Attila İlhan - "Ben Sana Mecburum" 
Ben sana mecburum bilemezsin
Adını mıh gibi aklımda tutuyorum
Büyüdükçe büyüyor gözlerin
Ben sana mecburum bilemezsin
İçimi seninle ısıtıyorum.'
$ws.Rows.Item(8).RowHeight = 86.4

$ws.Range("A9").Value = 'Cemal Süreya - "Tek Yasak"
Özgürlüğün geldiği gün,
O gün ölmek yasak!
Ve aşkın olduğu her yerde,
Yalnız kalmak yasak…'
$ws.Rows.Item(9).RowHeight = 72
$ws.Range("A9").WrapText = $true

$ws.Range("A10").Value = 'Yavuz Bülent Bakiler - "Söz"
Seninle bir yağmur başlıyor gözlerimde
Bir sevmek alıp başını gidiyor.
Öyle ak, öyle temiz ki bu sevgi
İnsan her şeyi seninle sevmek istiyor.'
$ws.Rows.Item(10).RowHeight = 86.4
$ws.Range("A10").WrapText = $true

$ws.Range("A11").Value = 'Cahit Külebi - "Hikâye" 
Senin dudakların pembe
Ellerin beyaz,
Al tut ellerimi bebek
Tut biraz!
Benim doğduğum köylerde
Ceviz ağaçları yoktu,
Bu yüzden serinliği sevmem,
Bu yüzden biraz hayalim.'
$ws.Rows.Item(11).RowHeight = 129.6
$ws.Range("A11").WrapText = $true

# Restore the selection shown in the saved workbook (cell A4).
$ws.Range("A4").Select() | Out-Null
